$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Nippon Large Cap Fund"
$ws.Range("B3").Value = 11110000
$ws.Range("C3").Value = "CMRPM0258F"
$ws.Range("D3").Value = "ASHAR "

$ws.Range("A4").Value = "Nippon Small Cap Fund"
$ws.Range("B4").Value = 1111
$ws.Range("C4").Value = "AGHPM9964E"
$ws.Range("D4").Value = "AMAN"

$ws.Range("D4").Select()
